$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change "Obrigatorio" column (E) from "N" to "S" for rows 2 through 11
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).Value = "S"
}

# Update Posicao_Inicio (B13) from 152 to 541
$ws.Cells.Item(13, 2).Value = 541
